$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above current row 2 ("50% overestimated"), shifting
# existing data rows (2-4) down to (3-5), to make room for the new
# "0% overestimated" scenario row.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the "0% overestimated" scenario data.
$ws.Range("A2").Value = "0% overestimated"
$ws.Range("B2").Value = 0.55405849026174403
$ws.Range("C2").Value = 0.621582133488505
$ws.Range("E2").Value = "0% overestimated"
$ws.Range("F2").Value = 0.69596385110244496
$ws.Range("G2").Value = 0.53694284724367503

# Update the active selection to match the authored workbook state.
$ws.Range("C2").Select()
